$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1686.2
$ws.Range("I17").Value = 1059
$ws.Range("J17").Value = 1999.8
$ws.Range("K17").Value = 3177
$ws.Range("L17").Value = 5999.4
$ws.Range("M17").Value = -3009
$ws.Range("N17").Value = -6335.4
$ws.Range("H19").Value = 191.42857
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 166.4
$ws.Range("K19").Value = 254
$ws.Range("L19").Value = 166.4
$ws.Range("M19").Value = -79
$ws.Range("N19").Value = -516.4
$ws.Range("H32").Value = 5000
$ws.Range("J32").Value = 5000
$ws.Range("L32").Value = 5000
$ws.Range("N32").Value = -5652
$ws.Range("H33").Value = 333.3846
$ws.Range("I33").Value = 118.7
$ws.Range("K33").Value = 118.7
$ws.Range("M33").Value = 110.3
$ws.Range("H34").Value = 4999.875
$ws.Range("I34").Value = 4857
$ws.Range("J34").Value = 6000
$ws.Range("K34").Value = 4857
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -4654
$ws.Range("N34").Value = -6406
$ws.Range("H36").Value = 4999.875
$ws.Range("I36").Value = 4857
$ws.Range("J36").Value = 6000
$ws.Range("K36").Value = 4857
$ws.Range("L36").Value = 6000
$ws.Range("M36").Value = -4142
$ws.Range("N36").Value = -7430
$ws.Range("H43").Value = 4006.5833
$ws.Range("I43").Value = 1618
$ws.Range("K43").Value = 1618
$ws.Range("M43").Value = -1549
$ws.Range("H55").Value = 458.3125
$ws.Range("I55").Value = 98.833336
$ws.Range("J55").Value = 674
$ws.Range("K55").Value = 98.833336
$ws.Range("L55").Value = 674
$ws.Range("M55").Value = 115.166664
$ws.Range("N55").Value = -1102
$ws.Range("H132").Value = 1529.1923
$ws.Range("I132").Value = 1531.1428
$ws.Range("K132").Value = 4593.428400000001
$ws.Range("M132").Value = -2063.428400000001
$ws.Range("H138").Value = 5437.52
$ws.Range("J138").Value = 5608.2827
$ws.Range("L138").Value = 16824.8481
$ws.Range("N138").Value = -27104.8481

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21000.814
$ws.Range("I32").Value = 11341.8
$ws.Range("J32").Value = 29399.957
$ws.Range("K32").Value = 11341.8
$ws.Range("L32").Value = 29399.957
$ws.Range("M32").Value = -11054.8
$ws.Range("N32").Value = -29973.957
$ws.Range("H122").Value = 591290
$ws.Range("I122").Value = 1001993
$ws.Range("J122").Value = 4571.4287
$ws.Range("K122").Value = 3005979
$ws.Range("L122").Value = 13714.2861
$ws.Range("M122").Value = -3003529
$ws.Range("N122").Value = -18614.2861
$ws.Range("H132").Value = 1853.7142
$ws.Range("I132").Value = 1765.8077
$ws.Range("J132").Value = 2996.5
$ws.Range("K132").Value = 5297.4231
$ws.Range("L132").Value = 8989.5
$ws.Range("M132").Value = -2767.4231
$ws.Range("N132").Value = -14049.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4109.3335
$ws.Range("I105").Value = 3356.8948
$ws.Range("K105").Value = 3356.8948
$ws.Range("M105").Value = -1609.8948
$ws.Range("H134").Value = 1252.92
$ws.Range("I134").Value = 692.3043
$ws.Range("J134").Value = 7700
$ws.Range("K134").Value = 2076.9129
$ws.Range("L134").Value = 23100
$ws.Range("M134").Value = 458.0870999999997
$ws.Range("N134").Value = -28170

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5561.533
$ws.Range("J31").Value = 6313.25
$ws.Range("L31").Value = 6313.25
$ws.Range("N31").Value = -6903.25
$ws.Range("H34").Value = 5561.533
$ws.Range("J34").Value = 6313.25
$ws.Range("L34").Value = 6313.25
$ws.Range("N34").Value = -6717.25
$ws.Range("H59").Value = 49999.25
$ws.Range("J59").Value = 49999.25
$ws.Range("L59").Value = 49999.25
$ws.Range("N59").Value = -52289.25
$ws.Range("H62").Value = 37086.582
$ws.Range("I62").Value = 3758.3333
$ws.Range("K62").Value = 3758.3333
$ws.Range("M62").Value = -3134.3333
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents() | Out-Null
$ws.Range("H65").Value = 37086.582
$ws.Range("I65").Value = 3758.3333
$ws.Range("K65").Value = 18791.6665
$ws.Range("M65").Value = -15671.6665
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents() | Out-Null
$ws.Range("H68").Value = 29623.5
$ws.Range("I68").Value = 24247.5
$ws.Range("K68").Value = 24247.5
$ws.Range("M68").Value = -23498.5
$ws.Range("H71").Value = 29623.5
$ws.Range("I71").Value = 24247.5
$ws.Range("K71").Value = 72742.5
$ws.Range("M71").Value = -68998.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 628.1429000000001
$ws.Range("J5").Value = 1405
$ws.Range("L5").Value = 4215
$ws.Range("N5").Value = -4439
$ws.Range("H12").Value = 13.888889
$ws.Range("I12").Value = 6
$ws.Range("J12").Value = 16.142857
$ws.Range("K12").Value = 18
$ws.Range("L12").Value = 48.428571
$ws.Range("M12").Value = 155
$ws.Range("N12").Value = -394.428571
$ws.Range("H26").Value = 484.8889
$ws.Range("I26").Value = 135
$ws.Range("J26").Value = 1184.6666
$ws.Range("K26").Value = 405
$ws.Range("L26").Value = 3553.9998
$ws.Range("M26").Value = -117
$ws.Range("N26").Value = -4129.9998
$ws.Range("H38").Value = 143.91667
$ws.Range("I38").Value = 61.5
$ws.Range("J38").Value = 308.75
$ws.Range("K38").Value = 184.5
$ws.Range("L38").Value = 926.25
$ws.Range("M38").Value = 162.5
$ws.Range("N38").Value = -1620.25
$ws.Range("H107").Value = 1062.5
$ws.Range("I107").Value = 592
$ws.Range("J107").Value = 1297.75
$ws.Range("K107").Value = 1776
$ws.Range("L107").Value = 3893.25
$ws.Range("M107").Value = 144
$ws.Range("N107").Value = -7733.25
$ws.Range("H135").Value = 628.1429000000001
$ws.Range("J135").Value = 1405
$ws.Range("L135").Value = 12645
$ws.Range("N135").Value = -17715

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3502.75
$ws.Range("I113").Value = 2255.5
$ws.Range("K113").Value = 2255.5
$ws.Range("M113").Value = -85.5
$ws.Range("H134").Value = 195000
$ws.Range("J134").Value = 195000
$ws.Range("L134").Value = 585000
$ws.Range("N134").Value = -590070

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2936.5
$ws.Range("I40").Value = 2936.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2936.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2800.5
$ws.Range("N40").ClearContents() | Out-Null
$ws.Range("H46").Value = 2950
$ws.Range("I46").Value = 2130
$ws.Range("K46").Value = 2130
$ws.Range("M46").Value = -1942
$ws.Range("H100").Value = 5628.2856
$ws.Range("I100").Value = 2399.6667
$ws.Range("J100").Value = 25000
$ws.Range("K100").Value = 2399.6667
$ws.Range("L100").Value = 25000
$ws.Range("M100").Value = -1858.6667
$ws.Range("N100").Value = -26082
$ws.Range("H122").Value = 2998
$ws.Range("I122").Value = 2998
$ws.Range("K122").Value = 8994
$ws.Range("M122").Value = -6544
$ws.Range("H132").Value = 5253.4614
$ws.Range("I132").Value = 5232.9165
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 15698.7495
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -13168.7495
$ws.Range("N132").Value = -21560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1200.3043
$ws.Range("I100").Value = 1468.875
$ws.Range("J100").Value = 586.4286
$ws.Range("K100").Value = 2937.75
$ws.Range("L100").Value = 1172.8572
$ws.Range("M100").Value = -2396.75
$ws.Range("N100").Value = -2254.8572
$ws.Range("H113").Value = 1834.1428
$ws.Range("I113").Value = 1622.8572
$ws.Range("K113").Value = 4868.571599999999
$ws.Range("M113").Value = -2698.571599999999
$ws.Range("H122").Value = 1926.75
$ws.Range("I122").Value = 1926.75
$ws.Range("K122").Value = 5780.25
$ws.Range("M122").Value = -3330.25
$ws.Range("H126").Value = 52027.6
$ws.Range("I126").Value = 67804.2
$ws.Range("K126").Value = 203412.6
$ws.Range("M126").Value = -200942.6
$ws.Range("H136").Value = 30642.176
$ws.Range("I136").Value = 939.6957
$ws.Range("K136").Value = 2819.0871
$ws.Range("M136").Value = -269.0870999999997
